# Update chart ideas - women
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Definitions")

# 1. "First Time Sailboard" row (row 6): update description text
$ws.Range("D6").Value = "All runs for top 10 first timers"

# 2. "Women's Fleet" row (row 10): rename, and strip the special yellow
#    highlight formatting (it becomes a normal row), drop its comment
$ws.Rows.Item(10).ClearFormats()
$ws.Range("B10:C10").HorizontalAlignment = -4108
$ws.Range("A10").Value = "Women's Fleet Sailboard"
$ws.Range("G10").Value = ""

# 3. Clear the "First Time Kite Fleet" row (row 11) - leave it blank
$ws.Range("A11:G11").ClearContents()

# 4. Remove the (now redundant) "First Time Wing Fleet" row entirely
$ws.Rows.Item(13).Delete()

# 5. Remove the (now redundant) "First Time Boat Fleet" row entirely
#    (it shifted up to row 14 after the previous delete)
$ws.Rows.Item(14).Delete()

# 6. Insert a new row above "Sailboards" (now at row 16) for the new
#    "Women" summary row, and give it the special highlighted style
$ws.Rows.Item(16).Insert()

$ws.Range("A16").Value = "Women"
$ws.Range("B16").Value = 1
$ws.Range("C16").Value = 10
$ws.Range("D16").Value = "All runs for top 10 women"
$ws.Range("E16").Value = "Different colour for each person. Light grey for non-top 10"
$ws.Range("F16").Value = "Larger marker for fastest run by each person"
$ws.Range("G16").Value = "No daily report but maybe it should be added?"

$ws.Range("A16:G16").Interior.Color = 65535
$ws.Range("B16:C16").HorizontalAlignment = -4108
